# Generate Report for Handback
# Updates the zh-cn and de-de handback status sheets:
#  - Status changes from "Ready for handoff" to "Handed back: in sync with en-US"
#  - Latest Target File / Latest Handback File columns (E/F) are now populated
#    with hyperlinks mirroring the Source File Name / Latest Handoff File links
#  - Latest Handback DateTime (G2) is stamped with the handback timestamp

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276  # BGR encoding of RGB FF6495ED (the workbook's HyperLink font color)

function Set-HandbackRow {
    param($SheetName, $MdAddress, $XlfAddress, $HandbackDateTime)

    $ws = $wb.Worksheets.Item($SheetName)

    # Status -> Handed back
    $ws.Range("B2").Value = "Handed back: in sync with en-US"

    # Latest Target File (E2): same file as the Source File Name hyperlink (A2)
    $ws.Hyperlinks.Add($ws.Range("E2"), $MdAddress, "", "", "076f1550-42b9-468a-9c2f-96707d73690b.md")
    $ws.Range("E2").Font.Underline = $true
    $ws.Range("E2").Font.Color = $hyperlinkColor

    # Latest Handback File (F2): same file as the Latest Handoff File hyperlink (C2)
    $ws.Hyperlinks.Add($ws.Range("F2"), $XlfAddress, "", "", $ws.Range("C2").Value2)
    $ws.Range("F2").Font.Underline = $true
    $ws.Range("F2").Font.Color = $hyperlinkColor

    # Latest Handback DateTime
    $ws.Range("G2").Value = $HandbackDateTime
}

Set-HandbackRow "zh-cn" `
    "https://github.com/OpenLocalizationTest/oltest/blob/002cf6e5da11df2979351765e2a7b79ca7eccff7/e2e/076f1550-42b9-468a-9c2f-96707d73690b.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/441fb6d49cfd1756d33f9e26f80a0cd480054086/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/076f1550-42b9-468a-9c2f-96707d73690b.1c4853cac13568b0bf9e5f983c7c7d3cc2fa3612.zh-cn.xlf" `
    "2016-01-28 04:56:35"

Set-HandbackRow "de-de" `
    "https://github.com/OpenLocalizationTest/oltest/blob/002cf6e5da11df2979351765e2a7b79ca7eccff7/e2e/076f1550-42b9-468a-9c2f-96707d73690b.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7e33946e4970a66f6d6fa69210bb82940124d73d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/076f1550-42b9-468a-9c2f-96707d73690b.1c4853cac13568b0bf9e5f983c7c7d3cc2fa3612.de-de.xlf" `
    "2016-01-28 04:56:53"

Write-Host "Handback report generated."
